# Auto-generated Excel COM-interop script to apply TPM data updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 14.11187666666667
$ws.Range("N2").Value = 42.33562999999999
$ws.Range("O2").Value = 0.08862966207485527
$ws.Range("P2").Value = 0.08862966207485526
$ws.Range("Q2").Value = 506.7148591205899
$ws.Range("R2").Value = 4560.433732085309
$ws.Range("S2").Value = 0.05412690701702674
$ws.Range("T2").Value = 0.05412690701702672

# Row 3
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.7176943460983047
$ws.Range("P3").Value = 0.7176943460983046
$ws.Range("Q3").Value = 4103.213088725299
$ws.Range("R3").Value = 36928.91779852769
$ws.Range("S3").Value = 0.4383021917098084
$ws.Range("T3").Value = 0.4383021917098082

# Row 4
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("O4").Value = 0.1936759918268401
$ws.Range("P4").Value = 0.1936759918268401
$ws.Range("Q4").Value = 1107.287341688063
$ws.Range("R4").Value = 9965.586075192567
$ws.Range("S4").Value = 0.1182796160521063
$ws.Range("T4").Value = 0.1182796160521062

# Row 5
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 14.11187666666667
$ws.Range("N5").Value = 42.33562999999999
$ws.Range("O5").Value = 0.08862966207485527
$ws.Range("P5").Value = 0.08862966207485526
$ws.Range("Q5").Value = 240.5880604085378
$ws.Range("R5").Value = 2165.29254367684
$ws.Range("S5").Value = 0.02569943892653964
$ws.Range("T5").Value = 0.02569943892653963

# Row 6
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.7176943460983047
$ws.Range("P6").Value = 0.7176943460983046
$ws.Range("S6").Value = 0.2081057468085388
$ws.Range("T6").Value = 0.2081057468085387

# Row 7
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("O7").Value = 0.1936759918268401
$ws.Range("P7").Value = 0.1936759918268401
$ws.Range("S7").Value = 0.05615912559033635
$ws.Range("T7").Value = 0.05615912559033633

# Row 8
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 14.11187666666667
$ws.Range("N8").Value = 42.33562999999999
$ws.Range("O8").Value = 0.08862966207485527
$ws.Range("P8").Value = 0.08862966207485526
$ws.Range("Q8").Value = 82.41319039081334
$ws.Range("R8").Value = 741.71871351732
$ws.Range("S8").Value = 0.008803316131288907
$ws.Range("T8").Value = 0.008803316131288904

# Row 9
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.7176943460983047
$ws.Range("P9").Value = 0.7176943460983046
$ws.Range("S9").Value = 0.07128640757995766
$ws.Range("T9").Value = 0.07128640757995765

# Row 10
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("O10").Value = 0.1936759918268401
$ws.Range("P10").Value = 0.1936759918268401
$ws.Range("S10").Value = 0.01923725018439753
$ws.Range("T10").Value = 0.01923725018439752
